$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Lrp6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01948966666666667
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.07096062449330311
$ws.Range("J2").Value = 0.07096062449330311
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.18802033333333
$ws.Range("N2").Value = 39.564061
$ws.Range("O2").Value = 0.1333092666742473
$ws.Range("P2").Value = 0.1413743409290003
$ws.Range("Q2").Value = 0.2570301202898889
$ws.Range("R2").Value = 2.313271082609
$ws.Range("S2").Value = 0.009459708813948867
$ws.Range("T2").Value = 0.010032011519651

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Lrp6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01948966666666667
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.07096062449330311
$ws.Range("J3").Value = 0.07096062449330311
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.845132
$ws.Range("N3").Value = 95.53539599999999
$ws.Range("O3").Value = 0.3219020813407859
$ws.Range("P3").Value = 0.3413768279472386
$ws.Range("Q3").Value = 0.6206510076359999
$ws.Range("R3").Value = 5.585859068724
$ws.Range("S3").Value = 0.02284237271763622
$ws.Range("T3").Value = 0.02422431289867894

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Lrp6"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01948966666666667
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.07096062449330311
$ws.Range("J4").Value = 0.07096062449330311
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 17.06250533333333
$ws.Range("N4").Value = 51.187516
$ws.Range("O4").Value = 0.172473958647377
$ws.Range("P4").Value = 0.1829084566999494
$ws.Range("Q4").Value = 0.3325425414448889
$ws.Range("R4").Value = 2.992882873004
$ws.Range("S4").Value = 0.01223885981445001
$ws.Range("T4").Value = 0.0129792983125347

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Lrp6"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01948966666666667
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.07096062449330311
$ws.Range("J5").Value = 0.07096062449330311
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.90151933333333
$ws.Range("N5").Value = 59.70455799999999
$ws.Range("O5").Value = 0.2011717362403738
$ws.Range("P5").Value = 0.2133424204787085
$ws.Range("Q5").Value = 0.3878739779668888
$ws.Range("R5").Value = 3.490865801701999
$ws.Range("S5").Value = 0.01427527203401898
$ws.Range("T5").Value = 0.01513891138808201

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Lrp6"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01948966666666667
$ws.Range("H6").Value = 0.058469
$ws.Range("I6").Value = 0.07096062449330311
$ws.Range("J6").Value = 0.07096062449330311
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.930832
$ws.Range("N6").Value = 33.861664
$ws.Range("O6").Value = 0.171142957097216
$ws.Range("P6").Value = 0.1209979539451033
$ws.Range("Q6").Value = 0.3299762720693333
$ws.Range("R6").Value = 1.979857632416
$ws.Range("S6").Value = 0.01214441111324903
$ws.Range("T6").Value = 0.008586090374356455

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Lrp6"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.255165
$ws.Range("H7").Value = 0.7654949999999999
$ws.Range("I7").Value = 0.9290393755066968
$ws.Range("J7").Value = 0.9290393755066969
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.18802033333333
$ws.Range("N7").Value = 39.564061
$ws.Range("O7").Value = 0.1333092666742473
$ws.Range("P7").Value = 0.1413743409290003
$ws.Range("Q7").Value = 3.365121208354999
$ws.Range("R7").Value = 30.28609087519499
$ws.Range("S7").Value = 0.1238495578602984
$ws.Range("T7").Value = 0.1313423294093493

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt1"
$ws.Range("C8").Value = "Lrp6"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.255165
$ws.Range("H8").Value = 0.7654949999999999
$ws.Range("I8").Value = 0.9290393755066968
$ws.Range("J8").Value = 0.9290393755066969
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 31.845132
$ws.Range("N8").Value = 95.53539599999999
$ws.Range("O8").Value = 0.3219020813407859
$ws.Range("P8").Value = 0.3413768279472386
$ws.Range("Q8").Value = 8.125763106779997
$ws.Range("R8").Value = 73.13186796101999
$ws.Range("S8").Value = 0.2990597086231497
$ws.Range("T8").Value = 0.3171525150485596

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt1"
$ws.Range("C9").Value = "Lrp6"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.255165
$ws.Range("H9").Value = 0.7654949999999999
$ws.Range("I9").Value = 0.9290393755066968
$ws.Range("J9").Value = 0.9290393755066969
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.06250533333333
$ws.Range("N9").Value = 51.187516
$ws.Range("O9").Value = 0.172473958647377
$ws.Range("P9").Value = 0.1829084566999494
$ws.Range("Q9").Value = 4.35375417338
$ws.Range("R9").Value = 39.18378756042
$ws.Range("S9").Value = 0.1602350988329269
$ws.Range("T9").Value = 0.1699291583874147

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt1"
$ws.Range("C10").Value = "Lrp6"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.255165
$ws.Range("H10").Value = 0.7654949999999999
$ws.Range("I10").Value = 0.9290393755066968
$ws.Range("J10").Value = 0.9290393755066969
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19.90151933333333
$ws.Range("N10").Value = 59.70455799999999
$ws.Range("O10").Value = 0.2011717362403738
$ws.Range("P10").Value = 0.2133424204787085
$ws.Range("Q10").Value = 5.078171180689998
$ws.Range("R10").Value = 45.70354062620999
$ws.Range("S10").Value = 0.1868964642063548
$ws.Range("T10").Value = 0.1982035090906265

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt1"
$ws.Range("C11").Value = "Lrp6"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.255165
$ws.Range("H11").Value = 0.7654949999999999
$ws.Range("I11").Value = 0.9290393755066968
$ws.Range("J11").Value = 0.9290393755066969
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 16.930832
$ws.Range("N11").Value = 33.861664
$ws.Range("O11").Value = 0.171142957097216
$ws.Range("P11").Value = 0.1209979539451033
$ws.Range("Q11").Value = 4.320155747279999
$ws.Range("R11").Value = 25.92093448368
$ws.Range("S11").Value = 0.158998545983967
$ws.Range("T11").Value = 0.1124118635707468
